# "calc toxicity at diff lengths"
#
# Adds a nontoxic-word counterpart to each existing toxic-word statistic:
#   - nontoxic_count  (new col D) = total_count - toxic_count
#   - nontoxic_pct    (new col G) = nontoxic_count / nontoxic_corpus_total * 100
# and renames the old difference / abs_difference columns (now I, J) to
# tox_total_diff / tox_total_abs_diff to reflect that they compare
# toxic_pct against total_pct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new columns; existing data shifts right and keeps its values,
# formatting and (for the header row) style.
#   before: A lemma B? ... (A index, B lemma, C toxic_count, D total_count, E toxic_pct, F total_pct, G difference, H abs_difference)
$ws.Columns("D:D").Insert()   # D becomes nontoxic_count; old D..H -> E..I
$ws.Columns("G:G").Insert()   # G becomes nontoxic_pct;   old E..I -> F..J

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, 4).Value  = "nontoxic_count"
$ws.Cells.Item(1, 7).Value  = "nontoxic_pct"
$ws.Cells.Item(1, 9).Value  = "tox_total_diff"
$ws.Cells.Item(1, 10).Value = "tox_total_abs_diff"

# --- Data rows ----------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Denominator for nontoxic_pct: total nontoxic-word occurrences across the
# full (Danish hate-speech) corpus this table was built from. This mirrors
# the pre-existing toxic_pct column, whose own (unlisted) denominator is the
# corpus's total toxic-word occurrence count (348) -- that figure isn't
# derivable from this identity-lemma subset either, it is a corpus-wide
# constant carried over from the original analysis.
$nontoxicCorpusTotal = 2283

for ($r = 2; $r -le $lastRow; $r++) {
    $toxicCount = $ws.Cells.Item($r, 3).Value()
    $totalCount = $ws.Cells.Item($r, 5).Value()
    $nontoxicCount = $totalCount - $toxicCount
    $ws.Cells.Item($r, 4).Value = $nontoxicCount
    $ws.Cells.Item($r, 7).Value = $nontoxicCount / $nontoxicCorpusTotal * 100
}
